$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update score values (Self-Evaluation Protocol score adjustments)
$ws.Range("C9").Value = 35
$ws.Range("C23").Value = 5
$ws.Range("C25").Value = 5
$ws.Range("C28").Value = 5

# Move the view/selection to match the new scroll position
$ws.Range("C26").Select()
